$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4429.2666
$ws.Range("J40").Value = 4791.4614
$ws.Range("L40").Value = 4791.4614
$ws.Range("N40").Value = -5141.4614
$ws.Range("H70").Value = 9505.786
$ws.Range("I70").Value = 5288.375
$ws.Range("J70").Value = 11192.75
$ws.Range("K70").Value = 15865.125
$ws.Range("L70").Value = 33578.25
$ws.Range("M70").Value = -15595.125
$ws.Range("N70").Value = -34118.25
$ws.Range("H73").Value = 9505.786
$ws.Range("I73").Value = 5288.375
$ws.Range("J73").Value = 11192.75
$ws.Range("K73").Value = 15865.125
$ws.Range("L73").Value = 33578.25
$ws.Range("M73").Value = -14929.125
$ws.Range("N73").Value = -35450.25
$ws.Range("H75").Value = 50285
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 50285
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H100").Value = 9314.056
$ws.Range("I100").Value = 3298.8333
$ws.Range("J100").Value = 12321.667
$ws.Range("K100").Value = 3298.8333
$ws.Range("L100").Value = 12321.667
$ws.Range("M100").Value = -2757.8333
$ws.Range("N100").Value = -13403.667
$ws.Range("H132").Value = 18869654
$ws.Range("I132").Value = 19609618
$ws.Range("K132").Value = 58828854
$ws.Range("M132").Value = -58826324
$ws.Range("H135").Value = 2280.4324
$ws.Range("I135").Value = 1142
$ws.Range("K135").Value = 10278
$ws.Range("M135").Value = -7743
$ws.Range("H137").Value = 2289.3953
$ws.Range("I137").Value = 2613.6667
$ws.Range("J137").Value = 1879.7894
$ws.Range("K137").Value = 7841.000100000001
$ws.Range("L137").Value = 5639.3682
$ws.Range("M137").Value = -5291.000100000001
$ws.Range("N137").Value = -10739.3682
$ws.Range("H138").Value = 2522.606
$ws.Range("I138").Value = 1039.2941
$ws.Range("J138").Value = 2830.122
$ws.Range("K138").Value = 3117.8823
$ws.Range("L138").Value = 8490.366
$ws.Range("M138").Value = 2022.1177
$ws.Range("N138").Value = -18770.366

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6227.4023
$ws.Range("I32").Value = 6067.256
$ws.Range("K32").Value = 6067.256
$ws.Range("M32").Value = -5780.256
$ws.Range("H45").Value = 2244.5293
$ws.Range("I45").Value = 2247.3125
$ws.Range("K45").Value = 2247.3125
$ws.Range("M45").Value = -1870.3125
$ws.Range("H74").Value = 4292.298
$ws.Range("I74").Value = 1434.9565
$ws.Range("J74").Value = 7030.5835
$ws.Range("K74").Value = 1434.9565
$ws.Range("L74").Value = 7030.5835
$ws.Range("M74").Value = -560.9565
$ws.Range("N74").Value = -8778.583500000001
$ws.Range("H77").Value = 4292.298
$ws.Range("I77").Value = 1434.9565
$ws.Range("J77").Value = 7030.5835
$ws.Range("K77").Value = 7174.7825
$ws.Range("L77").Value = 35152.9175
$ws.Range("M77").Value = -2806.7825
$ws.Range("N77").Value = -43888.9175
$ws.Range("H110").Value = 977.7778
$ws.Range("I110").Value = 280.4
$ws.Range("K110").Value = 280.4
$ws.Range("M110").Value = 1764.6
$ws.Range("H132").Value = 1415.8524
$ws.Range("I132").Value = 1231.12
$ws.Range("K132").Value = 3693.36
$ws.Range("M132").Value = -1163.36

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 4643
$ws.Range("I54").Value = 1915
$ws.Range("K54").Value = 1915
$ws.Range("M54").Value = -1431
$ws.Range("H64").Value = 1136.8125
$ws.Range("J64").Value = 1136.8125
$ws.Range("L64").Value = 1136.8125
$ws.Range("N64").Value = -1586.8125
$ws.Range("H67").Value = 1136.8125
$ws.Range("J67").Value = 1136.8125
$ws.Range("L67").Value = 1136.8125
$ws.Range("N67").Value = -2696.8125
$ws.Range("H86").Value = 3179859.2
$ws.Range("I86").Value = 4765136
$ws.Range("J86").Value = 9306.143
$ws.Range("K86").Value = 4765136
$ws.Range("L86").Value = 9306.143
$ws.Range("M86").Value = -4764013
$ws.Range("N86").Value = -11552.143
$ws.Range("H89").Value = 3179859.2
$ws.Range("I89").Value = 4765136
$ws.Range("J89").Value = 9306.143
$ws.Range("K89").Value = 23825680
$ws.Range("L89").Value = 46530.715
$ws.Range("M89").Value = -23820064
$ws.Range("N89").Value = -57762.715
$ws.Range("H107").Value = 4000
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 4000
$ws.Range("K107").Value = 0
$ws.Range("L107").ClearContents()
$ws.Range("M107").Value = 4000
$ws.Range("N107").Value = -7840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").ClearContents()
$ws.Range("N48").Value = 0
$ws.Range("H80").Value = 42248.75
$ws.Range("J80").Value = 42248.75
$ws.Range("L80").Value = 42248.75
$ws.Range("N80").Value = -44494.75
$ws.Range("H83").Value = 42248.75
$ws.Range("J83").Value = 42248.75
$ws.Range("L83").Value = 126746.25
$ws.Range("N83").Value = -137978.25
$ws.Range("H94").Value = 2077
$ws.Range("I94").Value = 1485.8
$ws.Range("J94").Value = 2372.6
$ws.Range("K94").Value = 1485.8
$ws.Range("L94").Value = 2372.6
$ws.Range("M94").Value = -1034.8
$ws.Range("N94").Value = -3274.6
$ws.Range("H99").Value = 4624.316
$ws.Range("I99").Value = 4772.636
$ws.Range("K99").Value = 4772.636
$ws.Range("M99").Value = -3274.636
$ws.Range("H126").Value = 4624.316
$ws.Range("I126").Value = 4772.636
$ws.Range("K126").Value = 14317.908
$ws.Range("M126").Value = -11847.908

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 200097
$ws.Range("J44").Value = 500094
$ws.Range("L44").Value = 1500282
$ws.Range("N44").Value = -1501078
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").ClearContents()
$ws.Range("N48").Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 30.733334
$ws.Range("I2").Value = 17.307692
$ws.Range("J2").Value = 118
$ws.Range("K2").Value = 17.307692
$ws.Range("L2").Value = 118
$ws.Range("M2").Value = 95.692308
$ws.Range("N2").Value = -344
$ws.Range("H39").Value = 94999.5
$ws.Range("J39").Value = 94999.5
$ws.Range("L39").Value = 94999.5
$ws.Range("N39").Value = -96063.5
$ws.Range("H80").Value = 3034.4614
$ws.Range("I80").Value = 2777.2144
$ws.Range("K80").Value = 2777.2144
$ws.Range("M80").Value = -1779.2144
$ws.Range("H83").Value = 3034.4614
$ws.Range("I83").Value = 2777.2144
$ws.Range("K83").Value = 13886.072
$ws.Range("M83").Value = -8894.072
$ws.Range("H102").Value = 22253.154
$ws.Range("I102").Value = 2640.5293
$ws.Range("J102").Value = 59299.223
$ws.Range("K102").Value = 2640.5293
$ws.Range("L102").Value = 59299.223
$ws.Range("M102").Value = -1018.5293
$ws.Range("N102").Value = -62543.223
$ws.Range("H132").Value = 5201.5713
$ws.Range("I132").Value = 5454.6924
$ws.Range("J132").Value = 1911
$ws.Range("K132").Value = 16364.0772
$ws.Range("L132").Value = 5733
$ws.Range("M132").Value = -13834.0772
$ws.Range("N132").Value = -10793

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 858.4595
$ws.Range("I16").Value = 361.03845
$ws.Range("K16").Value = 361.03845
$ws.Range("M16").Value = -191.03845
$ws.Range("H38").Value = 32515
$ws.Range("I38").Value = 30030
$ws.Range("J38").Value = 35000
$ws.Range("K38").Value = 30030
$ws.Range("L38").Value = 35000
$ws.Range("M38").Value = -29620
$ws.Range("N38").Value = -35820
$ws.Range("H82").Value = 5313.6924
$ws.Range("I82").Value = 7078.375
$ws.Range("K82").Value = 7078.375
$ws.Range("M82").Value = -6717.375
$ws.Range("H85").Value = 5313.6924
$ws.Range("I85").Value = 7078.375
$ws.Range("K85").Value = 7078.375
$ws.Range("M85").Value = -5830.375
$ws.Range("H132").Value = 3271.054
$ws.Range("J132").Value = 3267.138
$ws.Range("L132").Value = 9801.414000000001
$ws.Range("N132").Value = -14861.414

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 33333
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("H62").Value = 95882.89
$ws.Range("J62").Value = 13749
$ws.Range("L62").Value = 13749
$ws.Range("N62").Value = -14997
$ws.Range("H65").Value = 95882.89
$ws.Range("J65").Value = 13749
$ws.Range("L65").Value = 68745
$ws.Range("N65").Value = -74985
$ws.Range("H70").Value = 30000
$ws.Range("J70").Value = 30000
$ws.Range("L70").Value = 30000
$ws.Range("N70").Value = -30630
$ws.Range("H73").Value = 30000
$ws.Range("J73").Value = 30000
$ws.Range("L73").Value = 30000
$ws.Range("N73").Value = -32184
$ws.Range("H104").Value = 61372
$ws.Range("J104").Value = 61372
$ws.Range("L104").Value = 61372
$ws.Range("N104").Value = -68360
